$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# footer1.xml (COM: Footers.Item(2), wdHeaderFooterFirstPage) - PearsonLogo, docPr id="3"
$ftrFirst = $sec.Footers.Item(2)
$picF1 = $ftrFirst.Range.InlineShapes.Item(1)
$shpF1 = $picF1.ConvertToShape()
$shpF1.Name = "image2.png"
$shpF1.ConvertToInlineShape() | Out-Null

# footer2.xml (COM: Footers.Item(1), wdHeaderFooterPrimary) - PearsonLogo, docPr id="2"
$ftrPrimary = $sec.Footers.Item(1)
$picF2 = $ftrPrimary.Range.InlineShapes.Item(1)
$shpF2 = $picF2.ConvertToShape()
$shpF2.Name = "image2.png"
$shpF2.ConvertToInlineShape() | Out-Null

# header1.xml (COM: Headers.Item(2), wdHeaderFooterFirstPage) - BTec_Logo-Orange, docPr id="1"
$hdrFirst = $sec.Headers.Item(2)
$picH1 = $hdrFirst.Range.InlineShapes.Item(1)
$shpH1 = $picH1.ConvertToShape()
$shpH1.Name = "image1.jpg"
$shpH1.ConvertToInlineShape() | Out-Null
